$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.253.57"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "1.903.17"
$ws.Range("E3").Value = "  +1.37%  "
$ws.Range("E4").Value = "  -0.43%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.40%  "
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.692"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.42%  "
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.22"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.350"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0726"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.25%  "
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("D13").Value = "2.178.62"
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.92%  "
$ws.Range("E15").Value = "  +3.11%  "
$ws.Range("D16").Value = "1.910.32"
$ws.Range("E16").Value = "  +1.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "35.237.33"
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.86%  "
$ws.Range("D20").Value = "0.0₃0823"
$ws.Range("E20").Value = "  +2.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "241.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.01%  "
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("E25").Value = "  +1.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +14.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.48%  "
$ws.Range("E29").Value = "  +4.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "18.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.983"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.97%  "
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("E35").Value = "  -0.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.13"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.07%  "
$ws.Range("E38").Value = "  -1.03%  "
$ws.Range("E39").Value = "  -1.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0684"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +15.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("E42").Value = "  +2.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "90.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.60%  "
$ws.Range("D45").Value = "1.342.84"
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.11%  "
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("E50").Value = "  +1.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.19%  "

Write-Host "Applied all crypto list updates"
